$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove (clear) cells that no longer exist in the target state
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Update cells whose values changed (floating point precision corrections)
$ws.Range("E3").Value = 13.08276537368063

$ws.Range("C4").Value = -7.266312015249799
$ws.Range("E4").Value = 12.31225042954256

$ws.Range("C6").Value = 9.469137444079955
$ws.Range("E6").Value = 10.69920649119718

$ws.Range("C7").Value = 3.358206407534969
$ws.Range("E7").Value = 4.390489499870132

$ws.Range("E8").Value = 3.502435351035582

$ws.Range("C9").Value = 3.901355411819685
$ws.Range("E9").Value = 6.143002545701304

$ws.Range("E10").Value = 4.555278923792594

$ws.Range("E11").Value = 0.9515943257393467

$ws.Range("C12").Value = 5.246209615995689

$ws.Range("C13").Value = 4.862559663742938

$ws.Range("C14").Value = 2.76474001115945
$ws.Range("E14").Value = 0.3611963426345843

$ws.Range("C15").Value = -7.260793671746447
$ws.Range("E15").Value = 21.21858006100774

$ws.Range("C16").Value = 4.097586525396246

$ws.Range("C17").Value = 7.824284864703768

$ws.Range("C18").Value = -1.245022353133318
